# Auto-generated Excel COM-interop script
# Updates column F (want-to-go counts) across sheets
# per the source commit's regenerated data snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 990
$ws.Range("F4").Value = 8837
$ws.Range("F7").Value = 1951
$ws.Range("F8").Value = 6200
$ws.Range("F9").Value = 607
$ws.Range("F12").Value = 9151
$ws.Range("F13").Value = 10564
$ws.Range("F14").Value = 1203
$ws.Range("F15").Value = 1071
$ws.Range("F16").Value = 4813
$ws.Range("F17").Value = 757
$ws.Range("F18").Value = 405
$ws.Range("F20").Value = 316
$ws.Range("F22").Value = 1295
$ws.Range("F23").Value = 212
$ws.Range("F24").Value = 1851
$ws.Range("F26").Value = 1158
$ws.Range("F27").Value = 844
$ws.Range("F28").Value = 1988
$ws.Range("F30").Value = 575
$ws.Range("F31").Value = 2556
$ws.Range("F34").Value = 1642
$ws.Range("F35").Value = 88
$ws.Range("F37").Value = 402
$ws.Range("F38").Value = 888
$ws.Range("F40").Value = 3218
$ws.Range("F41").Value = 4204
$ws.Range("F42").Value = 228
$ws.Range("F43").Value = 78
$ws.Range("F44").Value = 481
$ws.Range("F45").Value = 556
$ws.Range("F46").Value = 23
$ws.Range("F48").Value = 223
$ws.Range("F49").Value = 4173

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F26").Value = 58

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5698

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 990
$ws.Range("F4").Value = 8837
$ws.Range("F8").Value = 6200
$ws.Range("F9").Value = 607
$ws.Range("F10").Value = 9151
$ws.Range("F11").Value = 9151
$ws.Range("F12").Value = 10564
$ws.Range("F14").Value = 1203
$ws.Range("F15").Value = 1071
$ws.Range("F16").Value = 4813
$ws.Range("F17").Value = 757
$ws.Range("F18").Value = 405
$ws.Range("F20").Value = 316
$ws.Range("F22").Value = 1295
$ws.Range("F23").Value = 212
$ws.Range("F25").Value = 1158
$ws.Range("F26").Value = 844
$ws.Range("F28").Value = 1988
$ws.Range("F30").Value = 575
$ws.Range("F31").Value = 2556
$ws.Range("F34").Value = 88
$ws.Range("F39").Value = 888
$ws.Range("F44").Value = 228
$ws.Range("F45").Value = 481
$ws.Range("F46").Value = 556
$ws.Range("F48").Value = 223
$ws.Range("F49").Value = 58

